$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9001
$ws.Range("I40").Value = 9001
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 9001
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8826

$ws.Range("H64").Value = 3200
$ws.Range("I64").Value = 2900
$ws.Range("K64").Value = 2900
$ws.Range("M64").Value = -2652

$ws.Range("H67").Value = 3200
$ws.Range("I67").Value = 2900
$ws.Range("K67").Value = 2900
$ws.Range("M67").Value = -2042

$ws.Range("H76").Value = 3010.7144
$ws.Range("I76").Value = 2696.3635
$ws.Range("J76").Value = 4163.3335
$ws.Range("K76").Value = 2696.3635
$ws.Range("L76").Value = 4163.3335
$ws.Range("M76").Value = -2381.3635
$ws.Range("N76").Value = -4793.3335

$ws.Range("H79").Value = 3010.7144
$ws.Range("I79").Value = 2696.3635
$ws.Range("J79").Value = 4163.3335
$ws.Range("K79").Value = 2696.3635
$ws.Range("L79").Value = 4163.3335
$ws.Range("M79").Value = -1604.3635
$ws.Range("N79").Value = -6347.3335

$ws.Range("H132").Value = 1429960.8
$ws.Range("I132").Value = 1931265.8
$ws.Range("J132").Value = 3169.4614
$ws.Range("K132").Value = 5793797.4
$ws.Range("L132").Value = 9508.3842
$ws.Range("M132").Value = -5791267.4
$ws.Range("N132").Value = -14568.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4522.42
$ws.Range("I32").Value = 3566.4043
$ws.Range("J32").Value = 19500
$ws.Range("K32").Value = 3566.4043
$ws.Range("L32").Value = 19500
$ws.Range("M32").Value = -3279.4043
$ws.Range("N32").Value = -20074

$ws.Range("H63").Value = 2355.2632
$ws.Range("I63").Value = 2330
$ws.Range("J63").Value = 2450
$ws.Range("K63").Value = 2330
$ws.Range("L63").Value = 2450
$ws.Range("M63").Value = -1644
$ws.Range("N63").Value = -3822

$ws.Range("H66").Value = 2355.2632
$ws.Range("I66").Value = 2330
$ws.Range("J66").Value = 2450
$ws.Range("K66").Value = 11650
$ws.Range("L66").Value = 12250
$ws.Range("M66").Value = -8218
$ws.Range("N66").Value = -19114

$ws.Range("H74").Value = 20834774
$ws.Range("I74").Value = 22728526
$ws.Range("J74").Value = 3507
$ws.Range("K74").Value = 22728526
$ws.Range("L74").Value = 3507
$ws.Range("M74").Value = -22727652
$ws.Range("N74").Value = -5255

$ws.Range("H77").Value = 20834774
$ws.Range("I77").Value = 22728526
$ws.Range("J77").Value = 3507
$ws.Range("K77").Value = 113642630
$ws.Range("L77").Value = 17535
$ws.Range("M77").Value = -113638262
$ws.Range("N77").Value = -26271

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2067945.6
$ws.Range("I105").Value = 3248471.8
$ws.Range("K105").Value = 3248471.8
$ws.Range("M105").Value = -3246724.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 251024.5
$ws.Range("I62").Value = 334000
$ws.Range("J62").Value = 2098
$ws.Range("K62").Value = 334000
$ws.Range("L62").Value = 2098
$ws.Range("M62").Value = -333376
$ws.Range("N62").Value = -3346

$ws.Range("H65").Value = 251024.5
$ws.Range("I65").Value = 334000
$ws.Range("J65").Value = 2098
$ws.Range("K65").Value = 1670000
$ws.Range("L65").Value = 10490
$ws.Range("M65").Value = -1666880
$ws.Range("N65").Value = -16730

$ws.Range("H107").Value = 1177.2727
$ws.Range("I107").Value = 609.25
$ws.Range("J107").Value = 1858.9
$ws.Range("K107").Value = 609.25
$ws.Range("L107").Value = 1858.9
$ws.Range("M107").Value = 1310.75
$ws.Range("N107").Value = -5698.9

$ws.Range("H132").Value = 1400.8163
$ws.Range("I132").Value = 1052.1025
$ws.Range("J132").Value = 2760.8
$ws.Range("K132").Value = 3156.3075
$ws.Range("L132").Value = 8282.400000000001
$ws.Range("M132").Value = -626.3074999999999
$ws.Range("N132").Value = -13342.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 603.4211
$ws.Range("I107").Value = 356.6
$ws.Range("J107").Value = 691.5714
$ws.Range("K107").Value = 1069.8
$ws.Range("L107").Value = 2074.7142
$ws.Range("M107").Value = 850.1999999999998
$ws.Range("N107").Value = -5914.7142

$ws.Range("H131").Value = 10913788
$ws.Range("J131").Value = 63476.188
$ws.Range("L131").Value = 190428.564
$ws.Range("N131").Value = -200508.564

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4936.857
$ws.Range("I70").Value = 5136
$ws.Range("J70").Value = 4787.5
$ws.Range("K70").Value = 5136
$ws.Range("L70").Value = 4787.5
$ws.Range("M70").Value = -4866
$ws.Range("N70").Value = -5327.5

$ws.Range("H73").Value = 4936.857
$ws.Range("I73").Value = 5136
$ws.Range("J73").Value = 4787.5
$ws.Range("K73").Value = 5136
$ws.Range("L73").Value = 4787.5
$ws.Range("M73").Value = -4200
$ws.Range("N73").Value = -6659.5

$ws.Range("H80").Value = 168967.5
$ws.Range("I80").Value = 2300
$ws.Range("J80").Value = 202301
$ws.Range("K80").Value = 2300
$ws.Range("L80").Value = 202301
$ws.Range("M80").Value = -1302
$ws.Range("N80").Value = -204297

$ws.Range("H83").Value = 168967.5
$ws.Range("I83").Value = 2300
$ws.Range("J83").Value = 202301
$ws.Range("K83").Value = 11500
$ws.Range("L83").Value = 1011505
$ws.Range("M83").Value = -6508
$ws.Range("N83").Value = -1021489

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 14578.223
$ws.Range("I68").Value = 23060.8
$ws.Range("J68").Value = 3975
$ws.Range("K68").Value = 23060.8
$ws.Range("L68").Value = 3975
$ws.Range("M68").Value = -22311.8
$ws.Range("N68").Value = -5473

$ws.Range("H71").Value = 14578.223
$ws.Range("I71").Value = 23060.8
$ws.Range("J71").Value = 3975
$ws.Range("K71").Value = 115304
$ws.Range("L71").Value = 19875
$ws.Range("M71").Value = -111560
$ws.Range("N71").Value = -27363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5084
$ws.Range("J45").Value = 5084
$ws.Range("L45").Value = 5084
$ws.Range("N45").Value = -6066

$ws.Range("H136").Value = 8094.722
$ws.Range("I136").Value = 2154.4443
$ws.Range("J136").Value = 14035
$ws.Range("K136").Value = 6463.3329
$ws.Range("L136").Value = 42105
$ws.Range("M136").Value = -3913.3329
$ws.Range("N136").Value = -47205
